$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 needs to pick up the header style (bold/border/centered) that the
# other header cells (B1:F1) already carry. Grab it via copy/paste-format
# while the original B1 cell (and its style) still exists.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Column A (rows 2:5) currently carries the bold/border style that the old,
# soon-to-be-dropped column A values had; the new column A (former column B)
# should be plain/default, so strip that formatting.
$ws.Range("A2:A5").ClearFormats()

# Capture current (pre-edit) values for columns A-F, rows 1-5 before we
# start overwriting them.
$orig = @{}
for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $orig["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# Shift everything one column to the left (old col A data is dropped
# entirely; old col B..F becomes new col A..E).
for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $orig["$r,$($c+1)"]
    }
}

# Drop the now-unused column F entirely.
$ws.Columns.Item(6).Delete() | Out-Null
